# ---------------------------------------------------------------------------
# "The Everlasting Symphony of Life" (Emily Dickinson)
#   -> "The Art of Capturing Time: An Introduction to Photography" (Lauren Giles)
#
# Rewrites the title/byline/contact line and the essay body/summary text,
# and appends the new "Body:" section with several new paragraphs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Plain substring replace that does NOT go through Word's typed-text
# autocorrect (so straight quotes/apostrophes are preserved verbatim).
function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                $true, 1, $false)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# Replace that goes through Find/Replace's "Replacement" text so that `^l`
# sequences turn into real line breaks (<w:br/>). Typed-text autocorrect can
# turn straight quotes into curly ones here, so avoid apostrophes in $new,
# or clean them up afterwards (see bottom of script).
function Replace-TextWithBreaks($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (breaks): $old"
    }
}

# --- Title / byline / contact -----------------------------------------------
Replace-Text "The Everlasting Symphony of Life" "The Art of Capturing Time: An Introduction to Photography"
Replace-Text "Emily Dickinson" "Lauren Giles"
Replace-Text "emilydickinson@poets" "laurengiles19@schoolmail"

# --- First body paragraph ----------------------------------------------------
Replace-Text "Amidst the vast tapestry of existence, life emerges as a symphony, an intricate harmony of natural phenomena" `
             "We live in a world of constant motion, where moments are fleeting and memories fade with time"

Replace-Text " From the tiniest of organisms to the grandest celestial bodies, all entities are interconnected, contributing to the enduring cacophony of life" `
             " Yet, there is a magical art form that has the power to freeze these fleeting instances, allowing us to revisit and relive them whenever we desire"

Replace-Text " In this symphony, each individual element, from the blooming of a flower to the migration of birds, plays a pivotal role, orchestrating a spectacle of beauty and wonder" `
             " This art form is photography"

Replace-Text "The intricate web of life begins with the fundamental building blocks - atoms and molecules - dancing in perpetual motion, forming inorganic matter" `
             "Photography is a unique blend of science and art"

Replace-Text " As the symphony of life evolves, atoms combine to form complex organic molecules, laying the groundwork for cellular life" `
             " It is rooted in the scientific understanding of light, optics, and chemistry, but it is also a creative expression that allows photographers to capture their unique perspectives and interpretations of the world"

# This also absorbs (and removes) the two runs that followed it in the
# original ("." and " Gradually, multicellular organisms ... symphony").
Replace-Text " The emergence of single-celled organisms marks a pivotal transition, as life begins to manifest its dazzling diversity. Gradually, multicellular organisms arise, forming intricate ecosystems, each species contributing its unique melody to the symphony" `
             " Through the lens of a camera, photographers can transform ordinary scenes into extraordinary works of art"

Replace-Text "From the depths of the oceans to the soaring heights of mountain ranges, life adapts and thrives in every conceivable environment" `
             "From the earliest days of photography in the mid-1800s, when pioneers like Nicephore Niepce and Louis Daguerre developed the first practical photographic processes, to the digital revolution that has made photography accessible to everyone, the art form has undergone remarkable changes"

Replace-Text " The Earth's diverse habitats, from lush rainforests to barren deserts, bear witness to the resilience and creativity of life" `
             " Yet, its essence remains the same: capturing moments in time and preserving them for posterity"

# The final sentence of the paragraph becomes just a line break, followed by
# the new "Body:" section (several new paragraphs-worth of text separated by
# manual line breaks, matching the source essay's style).
$bodySection = "^l" + `
    "^lBody:" + `
    "^l" + `
    "^lPhotography is an incredibly versatile and accessible art form." + `
    " With the advent of smartphones and digital cameras, anyone can take and share photographs, making it a truly democratic medium." + `
    " However, there is more to photography than simply pointing and shooting." + `
    " To create truly captivating images, photographers need to master the technical aspects of the craft, as well as develop their artistic vision." + `
    "^l" + `
    "^lOne of the key elements of photography is composition." + `
    " The photographer's choice of perspective, framing, and lighting can dramatically impact the mood and message of an image." + `
    " Whether it's the deliberate symmetry of a landscape photograph or the spontaneous chaos of a street scene, composition is the glue " + `
    "that holds a photograph together." + `
    "^l" + `
    "^lAnother important aspect of photography is editing." + `
    " In the digital age, photographers have access to a wide range of editing software that allows them to enhance and manipulate their images." + `
    " From simple adjustments like cropping and color correction to more complex techniques like compositing and HDR, editing can be used to bring out the best in a photograph and convey the photographer's intended message"

Replace-TextWithBreaks "Every living being, from microscopic bacteria to colossal whales, plays a crucial role in maintaining the delicate equilibrium of this global symphony, a tapestry of interconnectedness that sustains and nourishes all" `
             $bodySection

# --- Summary paragraph -------------------------------------------------------
Replace-Text "The symphony of life is an awe-inspiring testament to the intricate interconnectedness of all living beings" `
             "Photography is an art form that captures moments in time and preserves them for posterity"

Replace-Text " From the smallest molecules to the grandest ecosystems, each element contributes to the exquisite harmony of existence" `
             " It combines science, art, and creativity to create images that can inspire, inform, and connect people from all walks of life"

Replace-Text " The resilience and adaptability of life, as it thrives in diverse environments, exemplify the enduring power of nature's creative force" `
             " Through composition, editing, and a keen eye for detail, photographers can transform ordinary scenes into extraordinary works of art"

Replace-Text " This symphony reminds us that we are all integral parts of a larger collective, united by the common thread of life, inviting us to cherish and protect the magnificent tapestry of existence" `
             " Whether you're a seasoned photographer or just starting out, there is always something new to learn and explore in the world of photography"

# --- Trailing empty paragraph -------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null

# --- Clean-up: undo any smart-quote autocorrect from Replace-TextWithBreaks -
$quoteChars = @([char]0x2019, [char]0x2018, [char]0x201C, [char]0x201D)
foreach ($q in $quoteChars) {
    $straight = "'"
    if ($q -eq [char]0x201C -or $q -eq [char]0x201D) { $straight = '"' }
    $guard = 0
    $rng = $d.Content
    while ($rng.Find.Execute($q) -and $guard -lt 50) {
        $rng.Text = $straight
        $guard = $guard + 1
    }
}

Write-Output "edit complete"
